$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) and "全部类型" (sheet4) both contain the same table
# of convention listings and need the same numeric corrections applied.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("F2").Value = 860
    $ws.Range("F4").Value = 2170
    $ws.Range("F6").Value = 12662
    $ws.Range("F8").Value = 107
    $ws.Range("F10").Value = 463
    $ws.Range("F11").Value = 1152
    $ws.Range("F12").Value = 955
    $ws.Range("F13").Value = 13669
    $ws.Range("F14").Value = 14022
    $ws.Range("F23").Value = 1055

    # Row 26 item sold out: attendance count updated and price replaced
    # with a "sold out" label instead of a numeric price.
    $ws.Range("F26").Value = 928
    $ws.Range("G26").Value = "已售罄"

    $ws.Range("F27").Value = 5150
}
